$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2)
}

# Title (appears twice: H1 heading + bold run near the end)
Replace-Text "Play Braccio di Ferro Slot Game Free | RTP 96.67%" "Play Braccio di Ferro for Free at Top Casinos! (Review)"

# "What we like" bullet list
Replace-Text "Classic 5x3 layout with 30 fixed paylines" "Classic gameplay mechanics"
Replace-Text "Low minimum bet with preset automatic spins" "Low volatility and high RTP"
Replace-Text "Ample opportunities for big wins with a low volatility and high RTP" "Exciting bonus features"
Replace-Text "Unique bonus features including free spins and a Snakes and Ladders style game" "Captivating characters and theme"

# "What we don't like" bullet list
Replace-Text "Limited range of symbols" "Limited paylines (30 fixed paylines)"
Replace-Text "No progressive jackpot" "No progressive jackpot feature"

# Meta description italic run
Replace-Text "Read our review of Braccio di Ferro slot game and play for free. Enjoy ample opportunities for big wins with a low volatility and high RTP of 96.67%." "Play Braccio di Ferro for free and experience the excitement of this cult classic game with captivating characters and thrilling bonus features."
